# Generate Report for Handback
#
# The CI "Generate Report" pass ran again: a97c9924-dc0f-48ef-bd6e-6cacb517eb39
# came back from localization ("Handed back: in sync with en-US") while
# 9a7cd8ad-df34-4d55-b78d-e68cef751a0e is still only "Ready for handoff".
# Refresh the Overview sheet and both per-locale sheets to reflect that.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# URLs reused for the hyperlinks (same targets that already exist elsewhere
# in the workbook for these files).
# ---------------------------------------------------------------------------
$md_9a7cd8ad      = "https://github.com/OpenLocalizationTest/oltest/blob/d29181dcc378652fbf83055f69a7de986dd29221/e2e/9a7cd8ad-df34-4d55-b78d-e68cef751a0e.md"
$md_a97c9924      = "https://github.com/OpenLocalizationTest/oltest/blob/4bbeb0406562489dd8836c6b8ac1e6fcd306820d/e2e/a97c9924-dc0f-48ef-bd6e-6cacb517eb39.md"
$xlf_zhcn_9a7cd8ad = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7633c64d424c920535fef2cacb8521189e1a78f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9a7cd8ad-df34-4d55-b78d-e68cef751a0e.051aec48ec881805e40bf705fa8b90b6737dbbf0.zh-cn.xlf"
$xlf_zhcn_a97c9924 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/007c96866f8d4e0de4f2488e3f86949554ea45de/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a97c9924-dc0f-48ef-bd6e-6cacb517eb39.4463417102ea53953bebea6f1432c9e0b33ed4ed.zh-cn.xlf"
$xlf_dede_9a7cd8ad = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd3516129b983ce73121dbf58942f63e05060d4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9a7cd8ad-df34-4d55-b78d-e68cef751a0e.051aec48ec881805e40bf705fa8b90b6737dbbf0.de-de.xlf"
$xlf_dede_a97c9924 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75aec0556f0c6c9a073e6a2461660dffa0ec4a94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a97c9924-dc0f-48ef-bd6e-6cacb517eb39.4463417102ea53953bebea6f1432c9e0b33ed4ed.de-de.xlf"

$name_9a7cd8ad = "9a7cd8ad-df34-4d55-b78d-e68cef751a0e.md"
$name_a97c9924 = "a97c9924-dc0f-48ef-bd6e-6cacb517eb39.md"
$xlf_zhcn_9a7cd8ad_name = "9a7cd8ad-df34-4d55-b78d-e68cef751a0e.051aec48ec881805e40bf705fa8b90b6737dbbf0.zh-cn.xlf"
$xlf_zhcn_a97c9924_name = "a97c9924-dc0f-48ef-bd6e-6cacb517eb39.4463417102ea53953bebea6f1432c9e0b33ed4ed.zh-cn.xlf"
$xlf_dede_9a7cd8ad_name = "9a7cd8ad-df34-4d55-b78d-e68cef751a0e.051aec48ec881805e40bf705fa8b90b6737dbbf0.de-de.xlf"
$xlf_dede_a97c9924_name = "a97c9924-dc0f-48ef-bd6e-6cacb517eb39.4463417102ea53953bebea6f1432c9e0b33ed4ed.de-de.xlf"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady      = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: row 2 -> a97c9924 (handed back), row 3 -> 9a7cd8ad (ready)
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = $name_a97c9924
$ov.Range("B2").Value = $statusHandedBack
$ov.Range("C2").Value = $statusHandedBack
$ov.Range("D2").Value = "2016-29-18 14:29:58"

$ov.Range("A3").Value = $name_9a7cd8ad
$ov.Range("B3").Value = $statusReady
$ov.Range("C3").Value = $statusReady
$ov.Range("D3").Value = "2016-29-18 14:29:40"

$ov.Hyperlinks.Add($ov.Range("A2"), $md_a97c9924, "", "", $name_a97c9924) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), $md_9a7cd8ad, "", "", $name_9a7cd8ad) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: row 2 -> a97c9924 (handed back, now with target/handback info),
#              row 3 -> 9a7cd8ad (still just ready for handoff)
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $name_a97c9924
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $statusHandedBack
$zh.Range("D2").Value = $xlf_zhcn_a97c9924_name
$zh.Range("E2").Value = "2016-03-18 14:29:55"
$zh.Range("F2").Value = $name_a97c9924
$zh.Range("G2").Value = $xlf_zhcn_a97c9924_name
$zh.Range("H2").Value = "2016-03-18 14:30:22"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = $name_9a7cd8ad
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $statusReady
$zh.Range("D3").Value = $xlf_zhcn_9a7cd8ad_name
$zh.Range("E3").Value = "2016-03-18 14:29:37"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), $md_a97c9924, "", "", $name_a97c9924) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), $md_a97c9924, "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), $xlf_zhcn_a97c9924, "", "", $xlf_zhcn_a97c9924_name) | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), $md_a97c9924, "", "", $name_a97c9924) | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), $xlf_zhcn_a97c9924, "", "", $xlf_zhcn_a97c9924_name) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), $md_9a7cd8ad, "", "", $name_9a7cd8ad) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), $md_9a7cd8ad, "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), $xlf_zhcn_9a7cd8ad, "", "", $xlf_zhcn_9a7cd8ad_name) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same pattern as zh-cn, with de-de target/handback files
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = $name_a97c9924
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $statusHandedBack
$de.Range("D2").Value = $xlf_dede_a97c9924_name
$de.Range("E2").Value = "2016-03-18 14:29:58"
$de.Range("F2").Value = $name_a97c9924
$de.Range("G2").Value = $xlf_dede_a97c9924_name
$de.Range("H2").Value = "2016-03-18 14:30:28"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = $name_9a7cd8ad
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $statusReady
$de.Range("D3").Value = $xlf_dede_9a7cd8ad_name
$de.Range("E3").Value = "2016-03-18 14:29:40"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), $md_a97c9924, "", "", $name_a97c9924) | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), $md_a97c9924, "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), $xlf_dede_a97c9924, "", "", $xlf_dede_a97c9924_name) | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), $md_a97c9924, "", "", $name_a97c9924) | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), $xlf_dede_a97c9924, "", "", $xlf_dede_a97c9924_name) | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), $md_9a7cd8ad, "", "", $name_9a7cd8ad) | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), $md_9a7cd8ad, "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), $xlf_dede_9a7cd8ad, "", "", $xlf_dede_9a7cd8ad_name) | Out-Null

Write-Output "Generated handback report"
